$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "397.11" into numbers),
# matching the source workbook where these columns are stored as text.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '61.131.11'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '3.322.48'
$ws.Range("E3").Value = '  -0.64%  '

$ws.Range("E4").Value = '  +0.22%  '

Set-TextValue $ws.Range("D5") '397.11'
$ws.Range("E5").Value = '  -3.86%  '

Set-TextValue $ws.Range("D6") '124.03'
$ws.Range("E6").Value = '  +6.28%  '

Set-TextValue $ws.Range("D7") '0.583'
$ws.Range("E7").Value = '  +1.32%  '

$ws.Range("E8").Value = '  +0.03%  '

Set-TextValue $ws.Range("D9") '0.649'
$ws.Range("E9").Value = '  +2.87%  '

Set-TextValue $ws.Range("D10") '0.116'
$ws.Range("E10").Value = '  +0.32%  '

Set-TextValue $ws.Range("D11") '40.30'
$ws.Range("E11").Value = '  +0.09%  '

$ws.Range("E12").Value = '  -0.97%  '

$ws.Range("D13").Value = '3.853.45'
$ws.Range("E13").Value = '  -0.44%  '

Set-TextValue $ws.Range("D14") '8.12'
$ws.Range("E14").Value = '  -2.77%  '

Set-TextValue $ws.Range("D15") '18.99'
$ws.Range("E15").Value = '  -1.54%  '

$ws.Range("D16").Value = '3.410.68'
$ws.Range("E16").Value = '  +1.71%  '

$ws.Range("D17").Value = '61.083.13'
$ws.Range("E17").Value = '  -0.01%  '

Set-TextValue $ws.Range("D18") '11.03'
$ws.Range("E18").Value = '  +1.28%  '

Set-TextValue $ws.Range("D19") '0.991'
$ws.Range("E19").Value = '  -2.05%  '

Set-TextValue $ws.Range("D20") '0.0000125'
$ws.Range("E20").Value = '  +8.00%  '

Set-TextValue $ws.Range("D21") '3.15'
$ws.Range("E21").Value = '  -6.72%  '

Set-TextValue $ws.Range("D22") '78.81'
$ws.Range("E22").Value = '  +6.01%  '

Set-TextValue $ws.Range("D23") '12.54'
$ws.Range("E23").Value = '  -0.08%  '

Set-TextValue $ws.Range("D24") '295.37'
$ws.Range("E24").Value = '  -0.70%  '

Set-TextValue $ws.Range("D25") '3.06'
$ws.Range("E25").Value = '  -2.18%  '

$ws.Range("E26").Value = '  +11.43%  '

Set-TextValue $ws.Range("D27") '8.09'
$ws.Range("E27").Value = '  +6.92%  '

Set-TextValue $ws.Range("D28") '28.58'
$ws.Range("E28").Value = '  -2.48%  '

Set-TextValue $ws.Range("D29") '7.33'
$ws.Range("E29").Value = '  -6.28%  '

Set-TextValue $ws.Range("D30") '0.170'
$ws.Range("E30").Value = '  -0.78%  '

$ws.Range("E31").Value = '  -2.37%  '

$ws.Range("E32").Value = '  -0.10%  '

Set-TextValue $ws.Range("D33") '11.13'
$ws.Range("E33").Value = '  -1.85%  '

Set-TextValue $ws.Range("D34") '2.46'
$ws.Range("E34").Value = '  -2.88%  '

Set-TextValue $ws.Range("D35") '40.21'
$ws.Range("E35").Value = '  -5.99%  '

Set-TextValue $ws.Range("D36") '0.0468'
$ws.Range("E36").Value = '  -4.79%  '

Set-TextValue $ws.Range("D37") '51.87'
$ws.Range("E37").Value = '  -1.06%  '

Set-TextValue $ws.Range("D38") '0.999'
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("E39").Value = '  -2.79%  '

Set-TextValue $ws.Range("D40") '2.85'
$ws.Range("E40").Value = '  -7.52%  '

Set-TextValue $ws.Range("D41") '135.96'
$ws.Range("E41").Value = '  +0.73%  '

Set-TextValue $ws.Range("D42") '1.94'
$ws.Range("E42").Value = '  +1.49%  '

Set-TextValue $ws.Range("D43") '0.120'
$ws.Range("E43").Value = '  -0.54%  '

$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D44") '16.40'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D45") '0.274'
$ws.Range("E45").Value = '  -5.41%  '

Set-TextValue $ws.Range("D46") '3.78'
$ws.Range("E46").Value = '  -3.20%  '

$ws.Range("E47").Value = '  -0.85%  '

Set-TextValue $ws.Range("D48") '20.92'
$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '3.667.54'
$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.099.19'
$ws.Range("E50").Value = '  -2.59%  '

$ws.Range("E51").Value = '  -3.32%  '
